$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '30.531.87'
$ws.Range('E2').Value = '  +0.31%  '
$ws.Range('D3').Value = '2.139.08'
$ws.Range('E3').Value = '  +1.80%  '
$ws.Range('E4').Value = '  +0.61%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '352.41'
$ws.Range('E5').Value = '  +5.22%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.007'
$ws.Range('E6').Value = '  +0.47%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.5261'
$ws.Range('E7').Value = '  +0.68%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.4560'
$ws.Range('E8').Value = '  -0.03%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '53.48'
$ws.Range('E9').Value = '  -4.97%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.09163'
$ws.Range('E10').Value = '  +2.70%  '
$ws.Range('E11').Value = '  +1.38%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '25.55'
$ws.Range('E12').Value = '  +5.58%  '
$ws.Range('D13').Value = '2.136.31'
$ws.Range('E13').Value = '  +1.63%  '
$ws.Range('E14').Value = '  +1.33%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '8.202'
$ws.Range('E15').Value = '  +1.77%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '101.64'
$ws.Range('E16').Value = '  +4.73%  '
$ws.Range('E17').Value = '  +2.11%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '1.008'
$ws.Range('E18').Value = '  +0.28%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06725'
$ws.Range('E19').Value = '  +1.41%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '20.56'
$ws.Range('E20').Value = '  +7.12%  '
$ws.Range('E21').Value = '  +0.45%  '
$ws.Range('E22').Value = '  +1.34%  '
$ws.Range('D23').Value = '30.641.17'
$ws.Range('E23').Value = '  +0.49%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '12.87'
$ws.Range('E24').Value = '  +4.22%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.379'
$ws.Range('E25').Value = '  +0.83%  '
$ws.Range('D26').Value = '2.390.37'
$ws.Range('E26').Value = '  +1.79%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '22.55'
$ws.Range('E27').Value = '  +1.54%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.612'
$ws.Range('E28').Value = '  +3.88%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '165.07'
$ws.Range('E29').Value = '  +1.56%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '135.95'
$ws.Range('E30').Value = '  +2.09%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.226'
$ws.Range('E31').Value = '  +1.51%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.718'
$ws.Range('E32').Value = '  +3.63%  '
$ws.Range('E33').Value = '  +1.32%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '6.422'
$ws.Range('E34').Value = '  +0.58%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '4.037'
$ws.Range('E35').Value = '  +2.57%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '6.167'
$ws.Range('E36').Value = '  +4.66%  '
$ws.Range('E37').Value = '  +1.20%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.02649'
$ws.Range('E38').Value = '  +2.72%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.06995'
$ws.Range('E39').Value = '  +1.91%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.2360'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '12.73'
$ws.Range('E41').Value = '  +0.34%  '
$ws.Range('E42').Value = '  +1.72%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.278'
$ws.Range('E43').Value = '  +2.35%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '14.85'
$ws.Range('E44').Value = '  +6.64%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.356'
$ws.Range('E45').Value = '  +1.20%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.6521'
$ws.Range('E46').Value = '  +2.10%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.00000000376'
$ws.Range('E47').Value = '  +10.89%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '3.714'
$ws.Range('E48').Value = '  +1.53%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.251'
$ws.Range('E49').Value = '  +0.20%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '83.89'
$ws.Range('E50').Value = '  +0.84%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.07297'
$ws.Range('E51').Value = '  +2.46%  '
